$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above current row 5 (shifts old rows 5-23 down to 6-24)
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly record
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 45069
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101007
$ws.Range("J5").Value = "Kiwi"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 370
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19486
$ws.Range("Q5").Value = "`$/bandeja 18 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1083
$ws.Range("T5").Value = 18
